$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clarify "Price" column header -> "Basic Price"
$ws.Range("E1").Value = "Basic Price"

# Rename "GST Number" column header -> "GST Rate"
$ws.Range("F1").Value = "GST Rate"
